$wb = $excel.ActiveWorkbook

# --- Overall sheet: No. of 990 Filers w/ Gov Grants becomes formatted text ---
$wsOverall = $wb.Worksheets.Item("Overall")
$wsOverall.Range("A2").Value = "'3,264"

# --- County sheet: filer counts become text; append Total row 101 ---
$wsCounty = $wb.Worksheets.Item("County")
$wsCounty.Range("B2").Value = "'45"
$wsCounty.Range("B3").Value = "'9"
$wsCounty.Range("B4").Value = "'7"
$wsCounty.Range("B5").Value = "'2"
$wsCounty.Range("B6").Value = "'15"
$wsCounty.Range("B7").Value = "'8"
$wsCounty.Range("B8").Value = "'46"
$wsCounty.Range("B9").Value = "'8"
$wsCounty.Range("B10").Value = "'9"
$wsCounty.Range("B11").Value = "'28"
$wsCounty.Range("B12").Value = "'173"
$wsCounty.Range("B13").Value = "'34"
$wsCounty.Range("B14").Value = "'31"
$wsCounty.Range("B15").Value = "'26"
$wsCounty.Range("B16").Value = "'24"
$wsCounty.Range("B17").Value = "'7"
$wsCounty.Range("B18").Value = "'40"
$wsCounty.Range("B19").Value = "'28"
$wsCounty.Range("B20").Value = "'9"
$wsCounty.Range("B21").Value = "'5"
$wsCounty.Range("B22").Value = "'7"
$wsCounty.Range("B23").Value = "'32"
$wsCounty.Range("B24").Value = "'19"
$wsCounty.Range("B25").Value = "'26"
$wsCounty.Range("B26").Value = "'43"
$wsCounty.Range("B27").Value = "'3"
$wsCounty.Range("B28").Value = "'18"
$wsCounty.Range("B29").Value = "'37"
$wsCounty.Range("B30").Value = "'22"
$wsCounty.Range("B31").Value = "'18"
$wsCounty.Range("B32").Value = "'172"
$wsCounty.Range("B33").Value = "'15"
$wsCounty.Range("B34").Value = "'173"
$wsCounty.Range("B35").Value = "'11"
$wsCounty.Range("B36").Value = "'42"
$wsCounty.Range("B37").Value = "'1"
$wsCounty.Range("B38").Value = "'3"
$wsCounty.Range("B39").Value = "'18"
$wsCounty.Range("B40").Value = "'3"
$wsCounty.Range("B41").Value = "'167"
$wsCounty.Range("B42").Value = "'13"
$wsCounty.Range("B43").Value = "'23"
$wsCounty.Range("B44").Value = "'15"
$wsCounty.Range("B45").Value = "'48"
$wsCounty.Range("B46").Value = "'9"
$wsCounty.Range("B47").Value = "'11"
$wsCounty.Range("B48").Value = "'7"
$wsCounty.Range("B49").Value = "'30"
$wsCounty.Range("B50").Value = "'26"
$wsCounty.Range("B51").Value = "'37"
$wsCounty.Range("B52").Value = "'4"
$wsCounty.Range("B53").Value = "'20"
$wsCounty.Range("B54").Value = "'20"
$wsCounty.Range("B55").Value = "'12"
$wsCounty.Range("B56").Value = "'32"
$wsCounty.Range("B57").Value = "'15"
$wsCounty.Range("B58").Value = "'6"
$wsCounty.Range("B59").Value = "'16"
$wsCounty.Range("B60").Value = "'267"
$wsCounty.Range("B61").Value = "'13"
$wsCounty.Range("B62").Value = "'7"
$wsCounty.Range("B63").Value = "'39"
$wsCounty.Range("B64").Value = "'27"
$wsCounty.Range("B65").Value = "'75"
$wsCounty.Range("B66").Value = "'8"
$wsCounty.Range("B67").Value = "'32"
$wsCounty.Range("B68").Value = "'76"
$wsCounty.Range("B69").Value = "'9"
$wsCounty.Range("B70").Value = "'11"
$wsCounty.Range("B71").Value = "'14"
$wsCounty.Range("B72").Value = "'4"
$wsCounty.Range("B73").Value = "'7"
$wsCounty.Range("B74").Value = "'44"
$wsCounty.Range("B75").Value = "'15"
$wsCounty.Range("B76").Value = "'24"
$wsCounty.Range("B77").Value = "'14"
$wsCounty.Range("B78").Value = "'33"
$wsCounty.Range("B79").Value = "'16"
$wsCounty.Range("B80").Value = "'62"
$wsCounty.Range("B81").Value = "'45"
$wsCounty.Range("B82").Value = "'7"
$wsCounty.Range("B83").Value = "'10"
$wsCounty.Range("B84").Value = "'22"
$wsCounty.Range("B85").Value = "'12"
$wsCounty.Range("B86").Value = "'34"
$wsCounty.Range("B87").Value = "'7"
$wsCounty.Range("B88").Value = "'29"
$wsCounty.Range("B89").Value = "'5"
$wsCounty.Range("B90").Value = "'34"
$wsCounty.Range("B91").Value = "'18"
$wsCounty.Range("B92").Value = "'349"
$wsCounty.Range("B93").Value = "'10"
$wsCounty.Range("B94").Value = "'3"
$wsCounty.Range("B95").Value = "'40"
$wsCounty.Range("B96").Value = "'31"
$wsCounty.Range("B97").Value = "'23"
$wsCounty.Range("B98").Value = "'30"
$wsCounty.Range("B99").Value = "'15"
$wsCounty.Range("B100").Value = "'15"

# New Total row for County sheet (sheet grows from A1:F100 to A1:F101)
$wsCounty.Range("A101").Value = "Total"
$wsCounty.Range("B101").Value = "'3,264"
$wsCounty.Range("C101").Value = "'`$7,190,690,326"
$wsCounty.Range("D101").Value = "'9.29%"
$wsCounty.Range("E101").Value = "'-19.34%"
$wsCounty.Range("F101").Value = "'69.55%"

# --- Congressional District sheet: filer counts become text ---
$wsCD = $wb.Worksheets.Item("Congressional District")
$wsCD.Range("B2").Value = "'228"
$wsCD.Range("B3").Value = "'174"
$wsCD.Range("B4").Value = "'442"
$wsCD.Range("B5").Value = "'104"
$wsCD.Range("B6").Value = "'156"
$wsCD.Range("B7").Value = "'227"
$wsCD.Range("B8").Value = "'272"
$wsCD.Range("B9").Value = "'225"
$wsCD.Range("B10").Value = "'318"
$wsCD.Range("B11").Value = "'378"
$wsCD.Range("B12").Value = "'199"
$wsCD.Range("B13").Value = "'203"
$wsCD.Range("B14").Value = "'182"
$wsCD.Range("B15").Value = "'156"
$wsCD.Range("B16").Value = "'3,264"

# --- Size sheet: filer counts become text ---
$wsSize = $wb.Worksheets.Item("Size")
$wsSize.Range("B2").Value = "'1,199"
$wsSize.Range("B3").Value = "'839"
$wsSize.Range("B4").Value = "'584"
$wsSize.Range("B5").Value = "'221"
$wsSize.Range("B6").Value = "'288"
$wsSize.Range("B7").Value = "'133"
$wsSize.Range("B8").Value = "'3,264"

# --- Subsector sheet: filer counts become text ---
$wsSubsector = $wb.Worksheets.Item("Subsector")
$wsSubsector.Range("B2").Value = "'255"
$wsSubsector.Range("B3").Value = "'459"
$wsSubsector.Range("B4").Value = "'127"
$wsSubsector.Range("B5").Value = "'269"
$wsSubsector.Range("B6").Value = "'25"
$wsSubsector.Range("B7").Value = "'1,181"
$wsSubsector.Range("B8").Value = "'18"
$wsSubsector.Range("B9").Value = "'3"
$wsSubsector.Range("B10").Value = "'211"
$wsSubsector.Range("B11").Value = "'83"
$wsSubsector.Range("B12").Value = "'594"
$wsSubsector.Range("B13").Value = "'39"
$wsSubsector.Range("B14").Value = "'3,264"
